# Applies the "Lagt till Iterationsplan för Iteration 4" commit:
#  1. Updates the existing "Iteration 3" sheet (retrospective: statuses / actual
#     time filled in, one new row inserted for a task that was done).
#  2. Adds a brand-new sheet "Sheet1" (becomes the active tab) containing the
#     iteration plan for Iteration 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the existing "Iteration 3" worksheet
# ---------------------------------------------------------------------------
$it3 = $wb.Worksheets.Item("Iteration 3")

# Title now refers to "Iterationsplan Iteration 3" instead of "...Iteration 2"
$it3.Range("A1").Value = "Iterationsplan Iteration 3"

# Update status / actual-time columns for the already-existing rows
$it3.Cells.Item(8,3).Value = "Klar"
$it3.Cells.Item(8,5).Value = 1

$it3.Cells.Item(9,3).Value = "Klar"
$it3.Cells.Item(9,5).Value = 3

$it3.Cells.Item(10,3).Value = "Påbörjad"
$it3.Cells.Item(10,5).Value = 3

# Insert a brand-new row 11 for a task finished during the iteration
$it3.Rows.Item(11).Insert()
$it3.Cells.Item(11,1).Value = "F1"
$it3.Cells.Item(11,2).Value = "Gör så plusknappen skapar ny rad för ytterligare ingrediens"
$it3.Cells.Item(11,3).Value = "Påbörjad"
$it3.Cells.Item(11,4).Value = 1
$it3.Cells.Item(11,5).Value = 4

# Rows 12-14 (formerly 11-13) get status/actual time updates
$it3.Cells.Item(12,3).Value = "Påbörjad"
$it3.Cells.Item(12,5).Value = 1

$it3.Cells.Item(13,3).Value = "Klar"
$it3.Cells.Item(13,5).Value = 1

$it3.Cells.Item(14,3).Value = "Klar"
$it3.Cells.Item(14,4).Value = 1
$it3.Cells.Item(14,5).Value = 1

# Rows 15-16 (formerly 14-15) get actual-time values added
$it3.Cells.Item(15,5).Value = 0
$it3.Cells.Item(16,5).Value = 1

# Summary rows (formerly 16-18, now 17-19) get totals in column E
$it3.Cells.Item(17,5).Value = 15
$it3.Cells.Item(18,5).Value = 189
$it3.Cells.Item(19,5).Value = 240

# Sheet view: no longer the selected tab, selection becomes the title row
$it3.Range("A1:F1").Select()
$it3.Tab.Select()

# ---------------------------------------------------------------------------
# 2) Add a new "Sheet1" with the Iteration 4 plan (copy formatting from the
#    "Iteration 3" sheet, then overwrite contents)
# ---------------------------------------------------------------------------
$it3.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$it4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$it4.Name = "Sheet1"

# Remove the merged title cells present on the source sheet - the new sheet
# does not use merged cells for its header rows
$it4.Range("A1:F5").UnMerge()

# Header / title block
$it4.Range("A1").Value = "Iterationsplan Iteration 4"
$it4.Range("A2").Value = "Analys av föregående iteration"
$it4.Range("A3").Value = "Fixat så nya ingrediensrader dyker upp när man klickar på pluset. Även löst en bugg med min spinner/dropdownlista, men den funkar ännu inte i de nya ingrediensraderna."
$it4.Rows.Item(3).RowHeight = 31.5
$it4.Range("A4").Value = "Mål"
$it4.Range("A5").Value = "Målet med denna iteration är att kunna spara recepten lokalt och att kunna visa dem."

# Table header (row 7) stays the same as on "Iteration 3"

# Data rows 8-18
$it4.Cells.Item(8,1).Value = ""
$it4.Cells.Item(8,2).Value = "Handledarmöte"
$it4.Cells.Item(8,3).Value = "Klar"
$it4.Cells.Item(8,4).Value = 1
$it4.Cells.Item(8,5).Value = 1

$it4.Cells.Item(9,1).Value = "F1"
$it4.Cells.Item(9,2).Value = "Läs tutorial om att spara lokalt"
$it4.Cells.Item(9,3).Value = "Ej påbörjad"
$it4.Cells.Item(9,4).Value = 2
$it4.Cells.Item(9,5).Value = 2

$it4.Cells.Item(10,1).Value = "F1"
$it4.Cells.Item(10,2).Value = "Fixa så man kan spara lokalt"
$it4.Cells.Item(10,3).Value = "Ej påbörjad"
$it4.Cells.Item(10,4).Value = 8
$it4.Cells.Item(10,5).Value = 2

$it4.Cells.Item(11,1).Value = "F2"
$it4.Cells.Item(11,2).Value = "Skapa klickbar länk till sparat recept som visar receptet"
$it4.Cells.Item(11,3).Value = "Ej påbörjad"
$it4.Cells.Item(11,4).Value = 4
$it4.Cells.Item(11,5).Value = 1

$it4.Cells.Item(12,1).Value = "F1"
$it4.Cells.Item(12,2).Value = "Fixa så spinnern i nya ingrediensraderna fungerar"
$it4.Cells.Item(12,3).Value = "Ej påbörjad"
$it4.Cells.Item(12,4).Value = 3
$it4.Cells.Item(12,5).Value = ""

$it4.Cells.Item(13,1).Value = "F1"
$it4.Cells.Item(13,2).Value = "Gör så spinnerns värde skickar vidare och sparas"
$it4.Cells.Item(13,3).Value = "Ej påbörjad"
$it4.Cells.Item(13,4).Value = 2
$it4.Cells.Item(13,5).Value = ""

$it4.Cells.Item(14,1).Value = "F1"
$it4.Cells.Item(14,2).Value = "Gör så nya ingrediensradernas info också skickas vidare och sparas"
$it4.Cells.Item(14,3).Value = "Ej påbörjad"
$it4.Cells.Item(14,4).Value = 4
$it4.Cells.Item(14,5).Value = ""

$it4.Cells.Item(15,1).Value = "F4"
$it4.Cells.Item(15,2).Value = "Gör så man kan redigera ett existerande recept"
$it4.Cells.Item(15,3).Value = "Ej påbörjad"
$it4.Cells.Item(15,4).Value = 10
$it4.Cells.Item(15,5).Value = ""

$it4.Cells.Item(16,1).Value = "BK3"
$it4.Cells.Item(16,2).Value = "Fixa så texterna ""Ingrediens"" och ""Antal"" hamnar under nya ingrediensrader"
$it4.Cells.Item(16,3).Value = "Ej påbörjad"
$it4.Cells.Item(16,4).Value = 1
$it4.Cells.Item(16,5).Value = ""

# The copied sheet only has 19 rows (ending with the 3 summary rows at
# 17/18/19); insert two extra rows before the summary block for the two
# additional tasks, clearing the "summary row" style that Insert() carries
# over from the row it displaces.
$it4.Rows.Item(17).Insert()
$it4.Rows.Item(17).ClearFormats()
$it4.Rows.Item(17).RowHeight = $it4.Rows.Item(16).RowHeight
$it4.Rows.Item(17).Insert()
$it4.Rows.Item(17).ClearFormats()
$it4.Rows.Item(17).RowHeight = $it4.Rows.Item(16).RowHeight

$it4.Cells.Item(17,1).Value = ""
$it4.Cells.Item(17,2).Value = "Tester och testrapport"
$it4.Cells.Item(17,3).Value = "Ej påbörjad"
$it4.Cells.Item(17,4).Value = 4
$it4.Cells.Item(17,5).Value = ""

$it4.Cells.Item(18,1).Value = ""
$it4.Cells.Item(18,2).Value = "Skapa Iterationsplan för Iteration 5"
$it4.Cells.Item(18,3).Value = "Ej påbörjad"
$it4.Cells.Item(18,4).Value = 1
$it4.Cells.Item(18,5).Value = ""

# Summary rows 19-21 (formerly 17-19 on the copied sheet, now pushed down)
$it4.Cells.Item(19,3).Value = "Summa"
$it4.Cells.Item(19,4).Value = 40
$it4.Cells.Item(19,5).Value = ""

$it4.Cells.Item(20,3).Value = "Tid sedan föregående iteration"
$it4.Cells.Item(20,4).Value = ""
$it4.Cells.Item(20,5).Value = ""

$it4.Cells.Item(21,3).Value = "Total projekttid"
$it4.Cells.Item(21,4).Value = ""
$it4.Cells.Item(21,5).Value = 240

# Column widths specific to the new sheet
$it4.Columns.Item(2).ColumnWidth = 69.85546875
$it4.Columns.Item(6).ColumnWidth = 11.42578125

# Selection / active sheet
$it4.Range("F13").Select()
$it4.Activate()
